$wb = $excel.ActiveWorkbook

$oldGuid = "49352235-b7f5-4c20-a9f9-e09458ceaf52"
$newGuid = "b773c62c-bdbe-4c9b-9bcd-97003e97a4d2"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-23 11:01:40"

# Refresh the hyperlink display text on B2 (delete + re-add, same target URL)
$overviewLinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b9b5c95223e3c6563cd6bf9f161714ee38c1e181/e2e/$oldGuid.md"
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $overviewLinkUrl, "", "", "e2e\$newGuid.md")

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.18b4264ab6a0d1a86922fa65be54a341eada5835.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-23 11:01:35"
$wsZh.Range("I2").Value = ""
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"

# Drop the I2 hyperlink, keep A2's (new display text), same target URL
$zhLinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b9b5c95223e3c6563cd6bf9f161714ee38c1e181/e2e/$oldGuid.md"
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zhLinkUrl, "", "", "$newGuid.md")

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Range("G2").Value = "$newGuid.18b4264ab6a0d1a86922fa65be54a341eada5835.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-23 11:01:40"
$wsDe.Range("I2").Value = ""
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"

# Drop the I2 hyperlink, keep A2's (new display text), same target URL
$deLinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b9b5c95223e3c6563cd6bf9f161714ee38c1e181/e2e/$oldGuid.md"
$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $deLinkUrl, "", "", "$newGuid.md")
